$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (date serial, nuovi pos., somma mobile 7gg., somma mobile per 100mila ab.)
$data = @(
    @(230, 44304, 10, 48, 280.1283921797491),
    @(231, 44305, 11, 58, 338.4884738838634),
    @(232, 44306, 5, 51, 297.6364166909834),
    @(233, 44307, 6, 52, 303.4724248613948)
)

foreach ($row in $data) {
    $r = $row[0]

    # copy formatting from the last existing data row (229) so the new rows match style
    $ws.Range("A229:D229").Copy() | Out-Null
    $ws.Range("A$r`:D$r").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

$excel.CutCopyMode = 0
